# Add a new "Payment_Grade" worksheet with balance/monetization tiers and
# touch up the selection state on the sheets that were visited while doing
# so (mirrors the author's click-through before saving).

$wb = $excel.ActiveWorkbook

# --- 1. Create the new worksheet at the end of the workbook -----------------
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$ws = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $lastSheet)
$ws.Name = "Payment_Grade"

# --- 2. Header row ------------------------------------------------------
$ws.Range("A1").Value = "Grade"
$ws.Range("B1").Value = "Stat_Multiplier"
$ws.Range("C1").Value = "Note"

# --- 3. Data rows ---------------------------------------------------------
$ws.Range("A2").Value = "Free"
$ws.Range("B2").Value = 1
$ws.Range("C2").Value = "무과금 (기준)"

$ws.Range("A3").Value = "Light"
$ws.Range("B3").Value = 3
$ws.Range("C3").Value = "소과금 (월정액+패스)"

$ws.Range("A4").Value = "Heavy"
$ws.Range("B4").Value = 15
$ws.Range("C4").Value = "헤비과금 (랭커)"

# --- 4. Give the header row a thin box border (no bold/center, just border) -
$ws.Range("A1:C1").Borders.LineStyle = 1

# --- 5. Fit the columns to their new content --------------------------------
$ws.Range("A1:C4").EntireColumn.AutoFit() | Out-Null

# --- 6. Replay the selection / navigation state recorded in the workbook ----
$wb.Worksheets.Item("Growth_Table").Range("J27").Select() | Out-Null
$wb.Worksheets.Item("Skill_Data").Range("J28").Select() | Out-Null
$wb.Worksheets.Item("Monster_Book").Range("M29").Select() | Out-Null
$wb.Worksheets.Item("Dungeon_Config").Range("E8").Select() | Out-Null

# Land on the new sheet last so it becomes the active/selected tab.
$ws.Range("C5").Select() | Out-Null
